$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-26 down to 5-27.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new data.
$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value = "Bíobío"
$ws.Cells.Item(4, 4).Value = 44631
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 100112037
$ws.Cells.Item(4, 7).Value = "Cebollín"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 220
$ws.Cells.Item(4, 11).Value = 6000
$ws.Cells.Item(4, 12).Value = 6500
$ws.Cells.Item(4, 13).Value = 6227
$ws.Cells.Item(4, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(4, 15).Value = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value = 173
$ws.Cells.Item(4, 17).Value = 36
$ws.Cells.Item(4, 18).Value = "Hortaliza"
